$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.727484333333333
$ws.Range("H2").Value = 8.182453000000001
$ws.Range("I2").Value = 0.03096049453772388
$ws.Range("J2").Value = 0.03096049453772388
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 5.528925313526444
$ws.Range("R2").Value = 49.760327821738
$ws.Range("S2").Value = 0.0002042242322569154
$ws.Range("T2").Value = 0.0002042242322569154

$ws.Range("G3").Value = 2.727484333333333
$ws.Range("H3").Value = 8.182453000000001
$ws.Range("I3").Value = 0.03096049453772388
$ws.Range("J3").Value = 0.03096049453772388
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 699.4478251692165
$ws.Range("R3").Value = 6295.030426522948
$ws.Range("S3").Value = 0.025835797555356
$ws.Range("T3").Value = 0.025835797555356

$ws.Range("G4").Value = 2.727484333333333
$ws.Range("H4").Value = 8.182453000000001
$ws.Range("I4").Value = 0.03096049453772388
$ws.Range("J4").Value = 0.03096049453772388
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 133.2110594416712
$ws.Range("R4").Value = 1198.899534975041
$ws.Range("S4").Value = 0.004920472750110972
$ws.Range("T4").Value = 0.004920472750110971

$ws.Range("I5").Value = 0.5986009007423507
$ws.Range("J5").Value = 0.5986009007423507
$ws.Range("M5").Value = 2.027115333333333
$ws.Range("N5").Value = 6.081346
$ws.Range("O5").Value = 0.006596284565418616
$ws.Range("P5").Value = 0.006596284565418615
$ws.Range("Q5").Value = 106.8981527017116
$ws.Range("R5").Value = 962.0833743154039
$ws.Range("S5").Value = 0.003948541882412449
$ws.Range("T5").Value = 0.003948541882412449

$ws.Range("I6").Value = 0.5986009007423507
$ws.Range("J6").Value = 0.5986009007423507
$ws.Range("O6").Value = 0.8344762556643375
$ws.Range("P6").Value = 0.8344762556643374
$ws.Range("S6").Value = 0.4995182382887766
$ws.Range("T6").Value = 0.4995182382887765

$ws.Range("I7").Value = 0.5986009007423507
$ws.Range("J7").Value = 0.5986009007423507
$ws.Range("M7").Value = 48.84026566666667
$ws.Range("N7").Value = 146.520797
$ws.Range("O7").Value = 0.158927459770244
$ws.Range("P7").Value = 0.158927459770244
$ws.Range("Q7").Value = 2575.54865841912
$ws.Range("R7").Value = 23179.93792577208
$ws.Range("S7").Value = 0.09513412057116179
$ws.Range("T7").Value = 0.09513412057116175

$ws.Range("G8").Value = 32.63402300000001
$ws.Range("H8").Value = 97.90206900000001
$ws.Range("I8").Value = 0.3704386047199253
$ws.Range("J8").Value = 0.3704386047199253
$ws.Range("M8").Value = 2.027115333333333
$ws.Range("N8").Value = 6.081346
$ws.Range("O8").Value = 0.006596284565418616
$ws.Range("P8").Value = 0.006596284565418615
$ws.Range("Q8").Value = 66.15292841165267
$ws.Range("R8").Value = 595.3763557048741
$ws.Range("S8").Value = 0.002443518450749251
$ws.Range("T8").Value = 0.00244351845074925

$ws.Range("G9").Value = 32.63402300000001
$ws.Range("H9").Value = 97.90206900000001
$ws.Range("I9").Value = 0.3704386047199253
$ws.Range("J9").Value = 0.3704386047199253
$ws.Range("O9").Value = 0.8344762556643375
$ws.Range("P9").Value = 0.8344762556643374
$ws.Range("Q9").Value = 8368.809358467024
$ws.Range("R9").Value = 75319.28422620321
$ws.Range("S9").Value = 0.3091222198202048
$ws.Range("T9").Value = 0.3091222198202048

$ws.Range("G10").Value = 32.63402300000001
$ws.Range("H10").Value = 97.90206900000001
$ws.Range("I10").Value = 0.3704386047199253
$ws.Range("J10").Value = 0.3704386047199253
$ws.Range("M10").Value = 48.84026566666667
$ws.Range("N10").Value = 146.520797
$ws.Range("O10").Value = 0.158927459770244
$ws.Range("P10").Value = 0.158927459770244
$ws.Range("Q10").Value = 1593.854353092111
$ws.Range("R10").Value = 14344.689177829
$ws.Range("S10").Value = 0.05887286644897126
$ws.Range("T10").Value = 0.05887286644897124
